# Populate Sheet3 with the "Entity extract page" process notes, and make
# Sheet3 the active/selected sheet (moving tabSelected away from Sheet1).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws3 = $wb.Worksheets.Item(3)

# --- Sheet3 content -------------------------------------------------------
# Values are written in the same order the author originally typed them
# (this also controls the order new entries land in the shared-string
# table), then the numeric step markers are filled in.

$ws3.Range("D3").Value = "Greater than a certain length"
$ws3.Range("B3").Value = "Get blocks of significant text from the page"
$ws3.Range("C3").Value = "Significant"
$ws3.Range("D4").Value = "Not cookie or disclaimer"
$ws3.Range("D5").Value = "In a link but not a common use link"
$ws3.Range("D6").Value = "In a heading but not a common use heading"
$ws3.Range("B7").Value = "Glue together and tf-idf"
$ws3.Range("B9").Value = "Add weighting to highlighted words"
$ws3.Range("D9").Value = "Words in headings and significant links"
$ws3.Range("B13").Value = "Get links from page"
$ws3.Range("B14").Value = "Identifiy, about us, products/services/news"
$ws3.Range("B16").Value = "Identify link depth"
$ws3.Range("B15").Value = "Identify, meaningful domain specific links"
$ws3.Range("B17").Value = "Follow each link and repeat steps"
$ws3.Range("B20").Value = "Output"
$ws3.Range("C8").Value = "Tokenise words, noun phrases, names"
$ws3.Range("B1").Value = "Entity extract page"
$ws3.Range("B22").Value = "List of keywords and weightings"
$ws3.Range("C7").Value = "Remove stop words"

$ws3.Range("A3").Value = 1
$ws3.Range("A7").Value = 2
$ws3.Range("A9").Value = 3

# --- Column sizing ----------------------------------------------------
# Column B is best-fit to its (long) text entries; column C is a narrower
# manually-set width.

$ws3.Columns.Item(2).EntireColumn.AutoFit() | Out-Null
$ws3.Columns.Item(2).ColumnWidth = 40
$ws3.Columns.Item(3).ColumnWidth = 12.6

# --- Selection / active sheet ---------------------------------------------
# Activating Sheet3 and leaving the selection on L15 makes Sheet3 the
# workbook's active tab; Sheet1 automatically loses its tabSelected flag.

$ws3.Activate()
$ws3.Range("L15").Select()
